# Generate Report for Handback
# Fills in the "Latest Target File" (hyperlink), "Latest Handback File" and
# "Latest Handback DateTime" columns on the zh-cn / de-de sheets for both
# rows, flips the Overview sheet's per-language status text from
# "Ready for handoff" to "Handed back: in sync with en-US", and widens a
# few columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c65fc68726df8c4b9cf0daf0113bfc1fbf00aa9/e2e/"
$file1 = "2b6ee2ab-b090-4b82-a03a-4b2d2f0d859d.md"
$file2 = "3112559c-bec8-4e23-9d96-26a6970f9e3a.md"

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: status text for both languages, both rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# widen Overview columns E/F to fit the longer status text
$overview.Columns.Item(5).ColumnWidth = 29.14437166849777
$overview.Columns.Item(6).ColumnWidth = 29.14437166849777

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")

# "Status" column shared the same underlying text as the Overview sheet's
# per-language status cells, so it flips too.
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = "2b6ee2ab-b090-4b82-a03a-4b2d2f0d859d.aa5fac4323ee2e9aa90044c00ccfea98073466ae.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 08:30:52"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($urlBase + $file1), "", "", $file1)

$zhcn.Range("J3").Value = "3112559c-bec8-4e23-9d96-26a6970f9e3a.f2a1d8cd2a0871b75e621336f51ab0b0bd5c4a91.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 08:30:52"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($urlBase + $file2), "", "", $file2)

# widen columns to fit the new Status / Latest Target File / Latest Handback File text
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "2b6ee2ab-b090-4b82-a03a-4b2d2f0d859d.aa5fac4323ee2e9aa90044c00ccfea98073466ae.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 08:31:20"
$dede.Hyperlinks.Add($dede.Range("I2"), ($urlBase + $file1), "", "", $file1)

$dede.Range("J3").Value = "3112559c-bec8-4e23-9d96-26a6970f9e3a.f2a1d8cd2a0871b75e621336f51ab0b0bd5c4a91.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 08:31:20"
$dede.Hyperlinks.Add($dede.Range("I3"), ($urlBase + $file2), "", "", $file2)

# widen columns to fit the new Status / Latest Target File / Latest Handback File text
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report generated"
